# Weekly update for "Hortaliza, Vega Modelo de Temuco - Puerro":
# a new week's record is inserted above the existing row 280, pushing
# every following row (old 280-305) down by one and extending the used
# range from A1:R305 to A1:R306.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 280 (shifts rows 280:305 -> 281:306).
$ws.Rows.Item(280).Insert()

# Populate the newly inserted row with this week's data.
$ws.Cells.Item(280, 1).Value  = 10
$ws.Cells.Item(280, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(280, 3).Value  = "La Araucanía"
$ws.Cells.Item(280, 4).Value  = 45106
$ws.Cells.Item(280, 5).Value  = 9
$ws.Cells.Item(280, 6).Value  = 100112005
$ws.Cells.Item(280, 7).Value  = "Puerro"
$ws.Cells.Item(280, 8).Value  = "Azul de Maquehue"
$ws.Cells.Item(280, 9).Value  = "Primera"
$ws.Cells.Item(280, 10).Value = 110
$ws.Cells.Item(280, 11).Value = 10000
$ws.Cells.Item(280, 12).Value = 10000
$ws.Cells.Item(280, 13).Value = 10000
$ws.Cells.Item(280, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(280, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(280, 16).Value = 833
$ws.Cells.Item(280, 17).Value = 12
$ws.Cells.Item(280, 18).Value = "Hortaliza"
